$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("M3").Value = 1.07
$ws.Range("O3").Value = 1.07
$ws.Range("S3").Value = 1.87

# Row 4
$ws.Range("F4").Value = 2.38
$ws.Range("G4").Value = 2.8
$ws.Range("H4").Value = 3.55
$ws.Range("I4").Value = 4.7
$ws.Range("J4").Value = 2.46
$ws.Range("K4").Value = 3.1
$ws.Range("M4").Value = 1.14
$ws.Range("N4").Value = 2.2
$ws.Range("O4").Value = 1.66
$ws.Range("P4").Value = 1.39
$ws.Range("Q4").Value = 3
$ws.Range("S4").Value = 6
$ws.Range("T4").Value = 2.3
$ws.Range("U4").Value = 1.62
$ws.Range("V4").Value = 1.3
$ws.Range("W4").Value = 1.58
$ws.Range("X4").Value = 8.4
$ws.Range("Y4").Value = 9.6
$ws.Range("Z4").Value = 26
$ws.Range("AA4").Value = 120
$ws.Range("AB4").Value = 7
$ws.Range("AC4").Value = 7.4
$ws.Range("AD4").Value = 19
$ws.Range("AF4").Value = 15
$ws.Range("AG4").Value = 14
$ws.Range("AH4").Value = 29
$ws.Range("AJ4").Value = 44
$ws.Range("AL4").Value = 110
$ws.Range("AN4").Value = 65

# Row 5
$ws.Range("F5").Value = 1.41
$ws.Range("G5").Value = 1.5
$ws.Range("H5").Value = 9
$ws.Range("I5").Value = 13.5
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 4.8
$ws.Range("L5").Value = 1.46
$ws.Range("M5").Value = 1.09
$ws.Range("N5").Value = 2.88
$ws.Range("O5").Value = 1.43
$ws.Range("P5").Value = 1.63
$ws.Range("Q5").Value = 2.28
$ws.Range("R5").Value = 1.23
$ws.Range("S5").Value = 3.95
$ws.Range("U5").Value = 1.53
$ws.Range("V5").Value = 1.08
$ws.Range("W5").Value = 2.96
$ws.Range("X5").Value = 13.5
$ws.Range("Y5").Value = 30
$ws.Range("AB5").Value = 6.8
$ws.Range("AC5").Value = 13
$ws.Range("AD5").Value = 60
$ws.Range("AF5").Value = 8.199999999999999
$ws.Range("AG5").Value = 13.5
$ws.Range("AH5").Value = 48
$ws.Range("AJ5").Value = 14.5
$ws.Range("AK5").Value = 25
$ws.Range("AN5").Value = 13

# Row 6
$ws.Range("F6").Value = 2.32
$ws.Range("I6").Value = 3.25
$ws.Range("J6").Value = 3.65
$ws.Range("O6").Value = 1.27
$ws.Range("P6").Value = 2.12
$ws.Range("Q6").Value = 1.78
$ws.Range("R6").Value = 1.44
$ws.Range("S6").Value = 3
$ws.Range("U6").Value = 2.26
$ws.Range("V6").Value = 1.44
$ws.Range("AA6").Value = 55
$ws.Range("AO6").Value = 28

# Row 7
$ws.Range("F7").Value = 1.31
$ws.Range("G7").Value = 1.33
$ws.Range("H7").Value = 12
$ws.Range("J7").Value = 6.2
$ws.Range("K7").Value = 6.4
$ws.Range("S7").Value = 2.32
$ws.Range("T7").Value = 2.08
$ws.Range("U7").Value = 1.83
$ws.Range("V7").Value = 1.08
$ws.Range("W7").Value = 4
$ws.Range("Y7").Value = 48
$ws.Range("Z7").Value = 130
$ws.Range("AA7").Value = 510
$ws.Range("AB7").Value = 13.5
$ws.Range("AD7").Value = 46
$ws.Range("AE7").Value = 240
$ws.Range("AG7").Value = 14
$ws.Range("AH7").Value = 32
$ws.Range("AI7").Value = 170
$ws.Range("AJ7").Value = 10
$ws.Range("AK7").Value = 15
$ws.Range("AL7").Value = 1000
$ws.Range("AM7").Value = 170
$ws.Range("AN7").Value = 4.1
$ws.Range("AO7").Value = 320

# Row 8
$ws.Range("G8").Value = 1.68
$ws.Range("H8").Value = 6.2
$ws.Range("Q8").Value = 2.08
$ws.Range("T8").Value = 2.12
$ws.Range("U8").Value = 1.75
$ws.Range("W8").Value = 2.5
$ws.Range("AI8").Value = 120
$ws.Range("AJ8").Value = 15.5
$ws.Range("AL8").Value = 46
$ws.Range("AN8").Value = 1000

# Row 9
$ws.Range("F9").Value = 1.93
$ws.Range("N9").Value = 3.75
$ws.Range("T9").Value = 1.83
$ws.Range("U9").Value = 2.08

# Row 10
$ws.Range("L10").Value = 1.31
$ws.Range("S10").Value = 2.64

# Row 11
$ws.Range("F11").Value = 1.95
$ws.Range("G11").Value = 1.99
$ws.Range("I11").Value = 4.5
$ws.Range("K11").Value = 3.95
$ws.Range("V11").Value = 1.28
$ws.Range("W11").Value = 2
$ws.Range("AA11").Value = 110
$ws.Range("AD11").Value = 1000
$ws.Range("AE11").Value = 1000
$ws.Range("AG11").Value = 12

# Row 12
$ws.Range("Q12").Value = 2.04
$ws.Range("S12").Value = 3.75

# Row 13
$ws.Range("F13").Value = 3.25
$ws.Range("H13").Value = 2.64
$ws.Range("I13").Value = 2.88
$ws.Range("J13").Value = 2.84
$ws.Range("K13").Value = 2.88
$ws.Range("N13").Value = 1.31
$ws.Range("V13").Value = 1.53
$ws.Range("AF13").Value = 22
$ws.Range("AG13").Value = 17.5
$ws.Range("AK13").Value = 55
$ws.Range("AN13").Value = 70

# Row 14
$ws.Range("J14").Value = 7.4
$ws.Range("K14").Value = 8.199999999999999
$ws.Range("L14").Value = 1.28
$ws.Range("Q14").Value = 1.53
$ws.Range("S14").Value = 2.34
$ws.Range("X14").Value = 28

# Row 15
$ws.Range("H15").Value = 4.5
$ws.Range("N15").Value = 3.45
$ws.Range("P15").Value = 1.79
$ws.Range("R15").Value = 1.29
$ws.Range("S15").Value = 3.7

# Row 16
$ws.Range("J16").Value = 3.6
$ws.Range("N16").Value = 3.8
$ws.Range("P16").Value = 1.94
$ws.Range("AM16").Value = 95

# Row 17
$ws.Range("AF17").Value = 8.4

# Row 18
$ws.Range("I18").Value = 29
$ws.Range("T18").Value = 2.44

# Row 19
$ws.Range("AL19").Value = 46
$ws.Range("AN19").Value = 34
$ws.Range("AO19").Value = 16

# Row 20
$ws.Range("I20").Value = 13
$ws.Range("N20").Value = 7.2
$ws.Range("AE20").Value = 200
$ws.Range("AJ20").Value = 12

# Row 21
$ws.Range("N21").Value = 5.1
$ws.Range("P21").Value = 2.4
$ws.Range("S21").Value = 2.52
$ws.Range("T21").Value = 1.74
$ws.Range("AI21").Value = 1000
$ws.Range("AN21").Value = 6.8
$ws.Range("AO21").Value = 90

# Row 22
$ws.Range("K22").Value = 8.199999999999999
$ws.Range("L22").Value = 1.22
$ws.Range("N22").Value = 8.199999999999999
$ws.Range("R22").Value = 1.97
$ws.Range("S22").Value = 1.93

# Row 23
$ws.Range("F23").Value = 1.67
$ws.Range("H23").Value = 5.2
$ws.Range("I23").Value = 6
$ws.Range("K23").Value = 4.6
$ws.Range("S23").Value = 2.7
$ws.Range("W23").Value = 2.34
$ws.Range("AC23").Value = 1000
